# Apply the edit described in the commit:
# "Added aspects and classifications for inserting LCI data, fixed typos."
#
# 1. Fill in row 8 with a new user record (Paula Vollmer / pvollmer / IEF)
# 2. Rename the "Normal" cell style to "Standard" (typo fix)
# 3. Move the active selection to G11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New data row (row 8) -------------------------------------------------
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "pvollmer"
$ws.Range("C8").Value = "Paula Vollmer"
$ws.Range("D8").Value = "IEF"
$ws.Range("G8").Value = "2018-06-09 12:00:00"
$ws.Range("H8").Value = "2050-06-09 12:00:00"

# Match the date-column formatting used by the rows above (G/H use a
# quoted-text date style, not the plain placeholder style row 8 had before).
$ws.Range("G7:H7").Copy()
$ws.Range("G8:H8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fix typo in built-in cell style name ---------------------------------
$wb.Styles.Item("Normal").Name = "Standard"

# --- Move selection --------------------------------------------------------
$ws.Range("G11").Select()
